# Generate Report for Handback
#
# For the "zh-cn" and "de-de" status sheets, the two tracked files
# (rows 2 and 3) have come back from handback: in sync with en-US.
# Record that by:
#   - updating the Status text everywhere it is shown (Overview summary
#     sheet + each language sheet),
#   - filling in "Latest Target File" (E) / "Latest Handback File" (F)
#     with the same source/handoff file names (the handback round-tripped
#     cleanly), wired up as hyperlinks like the existing A/C columns,
#   - stamping "Latest Handback DateTime" (G) with the handback time.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$langSheets = @(
    @{ Name = "zh-cn"; G2Time = "2016-03-08 01:42:46"; G3Time = "2016-03-08 01:42:46";
       OrgSlug = "oltest.zh-cn"; CommitSha = "b911b212dcd17229a7aa03b11836a7e59c4d2eb4" },
    @{ Name = "de-de"; G2Time = "2016-03-08 01:43:04"; G3Time = "2016-03-08 01:43:04";
       OrgSlug = "oltest.de-de"; CommitSha = "67e47c9e467111d21c742a484cc9dfff4d230510" }
)

# --- Overview summary sheet: refresh the Status columns for both rows ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $statusText
$ov.Range("C2").Value = $statusText
$ov.Range("B3").Value = $statusText
$ov.Range("C3").Value = $statusText

foreach ($lang in $langSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Status column (B) for the two tracked files
    $ws.Range("B2").Value = $statusText
    $ws.Range("B3").Value = $statusText

    # Remember the existing handoff hyperlink targets (A2/C2/A3/C3/A4) and
    # their display text before wiping, so the whole collection can be
    # rebuilt in the right left-to-right, top-to-bottom order once the new
    # E/F handback columns are slotted in.
    $aDisplay2 = $ws.Range("A2").Value2
    $cDisplay2 = $ws.Range("C2").Value2
    $aDisplay3 = $ws.Range("A3").Value2
    $cDisplay3 = $ws.Range("C3").Value2
    $aDisplay4 = $ws.Range("A4").Value2

    $aTarget2 = "https://github.com/OpenLocalizationTest/oltest/blob/3c5a6e38570922dd1f6215a7732fb29b7eca66cc/e2e/$aDisplay2"
    $cTarget2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($lang.CommitSha)/ol-handoff/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/ht/$cDisplay2"
    $aTarget3 = "https://github.com/OpenLocalizationTest/oltest/blob/3c5a6e38570922dd1f6215a7732fb29b7eca66cc/e2e/$aDisplay3"
    $cTarget3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$($lang.CommitSha)/ol-handoff/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/ht/$cDisplay3"
    $aTarget4 = "https://github.com/OpenLocalizationTest/oltest/blob/3c5a6e38570922dd1f6215a7732fb29b7eca66cc/.localization-config"

    $eTarget2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.CommitSha)/ol-handback/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/hb/$aDisplay2"
    $fTarget2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.CommitSha)/ol-handback/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/hb/$cDisplay2"
    $eTarget3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.CommitSha)/ol-handback/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/hb/$aDisplay3"
    $fTarget3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$($lang.CommitSha)/ol-handback/OpenLocalizationTestOrg/$($lang.OrgSlug)/xinjiang/hb/$cDisplay3"

    # Fill in the handback mirror columns (values only for now; the
    # hyperlinks themselves are (re)built below in final left-to-right,
    # top-to-bottom order).
    $ws.Range("E2").Value = $aDisplay2
    $ws.Range("F2").Value = $cDisplay2
    $ws.Range("G2").Value = $lang.G2Time

    $ws.Range("E3").Value = $aDisplay3
    $ws.Range("F3").Value = $cDisplay3
    $ws.Range("G3").Value = $lang.G3Time

    # Wipe every hyperlink on the sheet and re-add them all in the final
    # order (A2, C2, E2, F2, A3, C3, E3, F3, A4) so the relationship ids
    # line up the way a fresh "generate handback report" pass would emit
    # them, instead of appending the new ones after the untouched A4 link.
    $ws.Range("A1:I4").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $aTarget2, "", "", $aDisplay2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $cTarget2, "", "", $cDisplay2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), $eTarget2, "", "", $aDisplay2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $fTarget2, "", "", $cDisplay2) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $aTarget3, "", "", $aDisplay3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $cTarget3, "", "", $cDisplay3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E3"), $eTarget3, "", "", $aDisplay3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $fTarget3, "", "", $cDisplay3) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $aTarget4, "", "", $aDisplay4) | Out-Null
}
